$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "=""54.582.41"""
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4163) | Out-Null
$ws.Range("E2").Value = "  -3.52%  "
$ws.Range("D3").Formula = "=""2.294.15"""
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4163) | Out-Null
$ws.Range("E3").Value = "  -3.27%  "
$ws.Range("D4").Formula = "=""1.00"""
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4163) | Out-Null
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Formula = "=""496.01"""
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Range("E5").Value = "  -2.42%  "
$ws.Range("D6").Formula = "=""127.45"""
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$ws.Range("E6").Value = "  -4.79%  "
$ws.Range("D7").Formula = "=""1.00"""
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  -2.37%  "
$ws.Range("D9").Formula = "=""2.293.71"""
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null
$ws.Range("E9").Value = "  -4.25%  "
$ws.Range("D10").Formula = "=""0.0951"""
$ws.Range("D10").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null
$ws.Range("E10").Value = "  -2.32%  "
$ws.Range("D11").Formula = "=""0.152"""
$ws.Range("D11").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163) | Out-Null
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Formula = "=""4.62"""
$ws.Range("D13").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163) | Out-Null
$ws.Range("E13").Value = "  -4.72%  "
$ws.Range("D14").Formula = "=""2.696.75"""
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4163) | Out-Null
$ws.Range("E14").Value = "  -3.76%  "
$ws.Range("D15").Formula = "=""21.64"""
$ws.Range("D15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("D16").Formula = "=""54.491.54"""
$ws.Range("D16").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4163) | Out-Null
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("E17").Value = "  -3.09%  "
$ws.Range("D18").Formula = "=""2.288.23"""
$ws.Range("D18").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163) | Out-Null
$ws.Range("E18").Value = "  -4.94%  "
$ws.Range("D19").Formula = "=""10.02"""
$ws.Range("D19").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D21").Formula = "=""304.17"""
$ws.Range("D21").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null
$ws.Range("E21").Value = "  -2.60%  "
$ws.Range("D22").Formula = "=""6.49"""
$ws.Range("D22").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Range("E22").Value = "  +3.15%  "
$ws.Range("D23").Formula = "=""0.999"""
$ws.Range("D23").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163) | Out-Null
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  -2.56%  "
$ws.Range("D25").Formula = "=""63.22"""
$ws.Range("D25").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null
$ws.Range("E25").Value = "  -3.80%  "
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("D29").Formula = "=""2.385.45"""
$ws.Range("D29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$ws.Range("E29").Value = "  -4.03%  "
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("D31").Formula = "=""170.64"""
$ws.Range("D31").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("E33").Value = "  -5.00%  "
$ws.Range("D34").Formula = "=""5.87"""
$ws.Range("D34").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163) | Out-Null
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("D35").Formula = "=""0.999"""
$ws.Range("D35").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163) | Out-Null
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("D37").Formula = "=""1.07"""
$ws.Range("D37").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$ws.Range("E37").Value = "  -2.81%  "
$ws.Range("D38").Formula = "=""17.58"""
$ws.Range("D38").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").Formula = "=""0.864"""
$ws.Range("D40").Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163) | Out-Null
$ws.Range("E40").Value = "  -1.93%  "
$ws.Range("E41").Value = "  -3.20%  "
$ws.Range("D42").Formula = "=""35.58"""
$ws.Range("D42").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$ws.Range("E42").Value = "  -2.91%  "
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("D45").Formula = "=""130.05"""
$ws.Range("D45").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("E47").Value = "  -4.07%  "
$ws.Range("D49").Formula = "=""0.548"""
$ws.Range("D49").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$ws.Range("E49").Value = "  -2.95%  "
$ws.Range("D50").Formula = "=""241.74"""
$ws.Range("D50").Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("E51").Value = "  -1.35%  "
$excel.CutCopyMode = $false
